$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Change 1: paragraph "Arpi (2012) nämner huvudtagligen My, en arbetare..."
#   - italic run "My" -> "Mys"
#   - split the trailing run so a new plain " " run is introduced
#   - rewrite the remainder of the sentence
# -------------------------------------------------------------------
$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("Arpi (2012) nämner huvudtagligen")) {
        $p1 = $cand
        break
    }
}

$myRng = $p1.Range
$found = $myRng.Find.Execute("My", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $myRng.Text = "Mys"

    $remStart = $myRng.End
    $remEnd = $p1.Range.End - 1
    $remRng = $d.Range($remStart, $remEnd)
    # Collapse the old trailing run down to a single plain space first, so
    # that the following insertion starts life as its own run (not merged
    # into the run that used to carry the comma sentence, nor inheriting
    # the italic formatting of "Mys").
    $remRng.Text = " "

    $afterSpace = $d.Range($remRng.End, $remRng.End)
    $afterSpace.InsertAfter("upplevelser i sin text. My är en arbetare på ett bemmaningsföretag som tog jobbet utav behov. Hennes position leder till att hon ibland får arbete och lönen som behövs, men ibland att hon inte får arbete alls.")
}

# -------------------------------------------------------------------
# Change 2: "Enligt arbettidslag förbjuds arbetsgivare att sätta..." paragraph
#   simple sentence-level rewrite (single run stays single run)
# -------------------------------------------------------------------
$old2 = "Enligt arbettidslag förbjuds arbetsgivare att sätta morgonspass på arbetare som arbetat kvällspass dagen innan berättar Arpi (2012); men enligt My så fick arbetare på vissa företag arbeta kvälls- och nästa dags morgonpass ändå."
$new2 = "Enligt arbettidslag förbjuds arbetsgivare från att sätta morgonspass på arbetare som arbetat kvällspass dagen innan berättar Arpi (2012); men enligt honom så påstår My att arbetare på vissa företag fick arbeta kvälls- och nästa dags morgonpass ändå."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# -------------------------------------------------------------------
# Change 3: merge the "Arpi (2012) inkluderar också statistik..." paragraph
# with the following "Slutligen så påstår Arpi (2012)..." paragraph into a
# single paragraph with new wording.
# -------------------------------------------------------------------
$statsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("Arpi (2012) inkluderar också statistik")) {
        $statsPara = $cand
        break
    }
}

if ($statsPara -ne $null) {
    # Delete the paragraph mark that separates it from the next paragraph,
    # merging the two paragraphs into one.
    $mark = $d.Range($statsPara.Range.End - 1, $statsPara.Range.End)
    $mark.Delete()

    $newText = "Slutligen så lyfter Arpi (2012) statistik ifrån SOU vilka pekar på att bemanningsbranschen tas upp huvudligen av personer där arbetet passar livssituationen. Dessutom så överrepresenteras vissa grupper, som unga, kvinnor, m.m. Arpi (2012) sammanfattar nämligen med att bemmaningsarbetare måste vara ytterst flexibla för att kunna klara av bemmanningsarbetet."
    $bodyRng = $d.Range($statsPara.Range.Start, $statsPara.Range.End - 1)
    $bodyRng.Text = $newText
}
